# It takes the number of rows from the training cases dinamically
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-case rows appended below the existing data (rows 5-15),
# reproducing the original 3 data rows (2-4) plus 4 more (5-8) and then
# the whole 7-row block duplicated again (9-15).
$data = @(
    @(5,  "Hombre", 22, "CABA",      3, "Universitario", 50000),
    @(6,  "Mujer",  21, "CABA",      5, "Posgrado",      60000),
    @(7,  "Otros",  20, "CABA",      7, "Posgrado",      70000),
    @(8,  "Otros",  20, "CABA",      7, "Posgrado",      70000),
    @(9,  "Hombre", 28, "Catamarca", 3, "Universitario", 80000),
    @(10, "Mujer",  30, "CABA",      5, "Posgrado",      120000),
    @(11, "Otros",  26, "GBA",       7, "Posgrado",      110000),
    @(12, "Hombre", 22, "CABA",      3, "Universitario", 50000),
    @(13, "Mujer",  21, "CABA",      5, "Posgrado",      60000),
    @(14, "Otros",  20, "CABA",      7, "Posgrado",      70000),
    @(15, "Otros",  20, "CABA",      7, "Posgrado",      70000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Widen column C so the newly inserted text values fit.
# (COM ColumnWidth and the stored OOXML width differ by the default
# font's padding offset; 14.1666... round-trips to a stored width of 15.)
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666

# Match the recorded selection after the edit.
$ws.Range("H11").Select()
